$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 3750184.5
$ws.Range("J4").Value = 20000332
$ws.Range("L4").Value = 20000332
$ws.Range("N4").Value = -20000560
$ws.Range("H9").Value = 417.33334
$ws.Range("I9").Value = 99
$ws.Range("J9").Value = 429.57693
$ws.Range("K9").Value = 99
$ws.Range("L9").Value = 429.57693
$ws.Range("M9").Value = 70
$ws.Range("N9").Value = -767.5769299999999
$ws.Range("H11").Value = 54.541668
$ws.Range("I11").Value = 54.541668
$ws.Range("K11").Value = 54.541668
$ws.Range("M11").Value = 85.458332
$ws.Range("H17").Value = 1043769.4
$ws.Range("I17").Value = 1257.3334
$ws.Range("J17").Value = 1490560.2
$ws.Range("K17").Value = 3772.0002
$ws.Range("L17").Value = 4471680.6
$ws.Range("M17").Value = -3604.0002
$ws.Range("N17").Value = -4472016.6
$ws.Range("H28").Value = 1269.1666
$ws.Range("I28").Value = 1302.25
$ws.Range("K28").Value = 1302.25
$ws.Range("M28").Value = -817.25
$ws.Range("H100").Value = 4171.625
$ws.Range("I100").Value = 2114
$ws.Range("K100").Value = 2114
$ws.Range("M100").Value = -1573
$ws.Range("H107").Value = 1619.1875
$ws.Range("I107").Value = 1619.1875
$ws.Range("K107").Value = 1619.1875
$ws.Range("M107").Value = 300.8125
$ws.Range("H113").Value = 3803.1
$ws.Range("I113").Value = 3753.6924
$ws.Range("J113").Value = 4124.25
$ws.Range("K113").Value = 3753.6924
$ws.Range("L113").Value = 4124.25
$ws.Range("M113").Value = -499.6923999999999
$ws.Range("N113").Value = -10632.25
$ws.Range("H129").Value = 1480.9286
$ws.Range("I129").Value = 856.8889
$ws.Range("J129").Value = 2604.2
$ws.Range("K129").Value = 2570.6667
$ws.Range("L129").Value = 7812.599999999999
$ws.Range("M129").Value = 2429.3333
$ws.Range("N129").Value = -17812.6
$ws.Range("H141").Value = 5722.387
$ws.Range("I141").Value = 3269.4583
$ws.Range("J141").Value = 14132.429
$ws.Range("K141").Value = 9808.374899999999
$ws.Range("L141").Value = 42397.287
$ws.Range("M141").Value = -4628.374899999999
$ws.Range("N141").Value = -52757.287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2752.439
$ws.Range("I74").Value = 2566.923
$ws.Range("J74").Value = 3074
$ws.Range("K74").Value = 2566.923
$ws.Range("L74").Value = 3074
$ws.Range("M74").Value = -1692.923
$ws.Range("N74").Value = -4822
$ws.Range("H77").Value = 2752.439
$ws.Range("I77").Value = 2566.923
$ws.Range("J77").Value = 3074
$ws.Range("K77").Value = 12834.615
$ws.Range("L77").Value = 15370
$ws.Range("M77").Value = -8466.614999999998
$ws.Range("N77").Value = -24106
$ws.Range("H110").Value = 1080
$ws.Range("I110").Value = 1045
$ws.Range("J110").Value = 1500
$ws.Range("K110").Value = 1045
$ws.Range("L110").Value = 1500
$ws.Range("M110").Value = 1000
$ws.Range("N110").Value = -5590

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3105.739
$ws.Range("I99").Value = 2345.0625
$ws.Range("J99").Value = 4844.4287
$ws.Range("K99").Value = 2345.0625
$ws.Range("L99").Value = 4844.4287
$ws.Range("M99").Value = -847.0625
$ws.Range("N99").Value = -7840.4287
$ws.Range("H105").Value = 2583.075
$ws.Range("I105").Value = 2514.92
$ws.Range("K105").Value = 2514.92
$ws.Range("M105").Value = -767.9200000000001
$ws.Range("H107").Value = 3451.2888
$ws.Range("I107").Value = 2359.9443
$ws.Range("K107").Value = 2359.9443
$ws.Range("M107").Value = -439.9443000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2373.0544
$ws.Range("I31").Value = 1590.5454
$ws.Range("K31").Value = 1590.5454
$ws.Range("M31").Value = -1295.5454
$ws.Range("H34").Value = 2373.0544
$ws.Range("I34").Value = 1590.5454
$ws.Range("K34").Value = 1590.5454
$ws.Range("M34").Value = -1388.5454
$ws.Range("H99").Value = 2386.923
$ws.Range("J99").Value = 2456.1428
$ws.Range("L99").Value = 2456.1428
$ws.Range("N99").Value = -5452.1428
$ws.Range("H105").Value = 3903.2778
$ws.Range("I105").Value = 5238.625
$ws.Range("K105").Value = 5238.625
$ws.Range("M105").Value = -3491.625
$ws.Range("H107").Value = 822.06665
$ws.Range("J107").Value = 672
$ws.Range("L107").Value = 672
$ws.Range("N107").Value = -4512
$ws.Range("H126").Value = 2386.923
$ws.Range("J126").Value = 2456.1428
$ws.Range("L126").Value = 7368.428400000001
$ws.Range("N126").Value = -12308.4284
$ws.Range("H132").Value = 1830
$ws.Range("I132").Value = 1063
$ws.Range("K132").Value = 3189
$ws.Range("M132").Value = -659
$ws.Range("H141").Value = 268045.7
$ws.Range("J141").Value = 268045.7
$ws.Range("L141").Value = 268045.7
$ws.Range("N141").Value = -278405.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 39009628
$ws.Range("I4").Value = 40489564
$ws.Range("K4").Value = 121468692
$ws.Range("M4").Value = -121468580
$ws.Range("H11").Value = 60362
$ws.Range("I11").Value = 60362
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 181086
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -180946
$ws.Range("N11").ClearContents()
$ws.Range("H57").Value = 4549.375
$ws.Range("I57").Value = 465
$ws.Range("J57").Value = 7000
$ws.Range("K57").Value = 1395
$ws.Range("L57").Value = 21000
$ws.Range("M57").Value = -836
$ws.Range("N57").Value = -22118
$ws.Range("H140").Value = 1788.7826
$ws.Range("I140").Value = 1724.6364
$ws.Range("K140").Value = 5173.9092
$ws.Range("M140").Value = 6.090799999999945

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4274.8125
$ws.Range("I7").Value = 3628.4
$ws.Range("J7").Value = 5352.1665
$ws.Range("K7").Value = 3628.4
$ws.Range("L7").Value = 5352.1665
$ws.Range("M7").Value = -3516.4
$ws.Range("N7").Value = -5576.1665
$ws.Range("H40").Value = 5153.433
$ws.Range("I40").Value = 4211
$ws.Range("J40").Value = 6567.0835
$ws.Range("K40").Value = 4211
$ws.Range("L40").Value = 6567.0835
$ws.Range("M40").Value = -4075
$ws.Range("N40").Value = -6839.0835
$ws.Range("H126").Value = 4274.8125
$ws.Range("I126").Value = 3628.4
$ws.Range("J126").Value = 5352.1665
$ws.Range("K126").Value = 10885.2
$ws.Range("L126").Value = 16056.4995
$ws.Range("M126").Value = -8415.200000000001
$ws.Range("N126").Value = -20996.4995
$ws.Range("H132").Value = 4257.115
$ws.Range("I132").Value = 2423.1177
$ws.Range("J132").Value = 7721.3335
$ws.Range("K132").Value = 7269.353099999999
$ws.Range("L132").Value = 23164.0005
$ws.Range("M132").Value = -4739.353099999999
$ws.Range("N132").Value = -28224.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 84999.75
$ws.Range("J94").Value = 84999.75
$ws.Range("L94").Value = 84999.75
$ws.Range("N94").Value = -86801.75
$ws.Range("H107").Value = 905.8
$ws.Range("I107").Value = 944.75
$ws.Range("K107").Value = 2834.25
$ws.Range("M107").Value = -914.25
$ws.Range("H117").Value = 52499
$ws.Range("J117").Value = 52499
$ws.Range("L117").Value = 52499
$ws.Range("N117").Value = -61677
$ws.Range("H122").Value = 4588.591
$ws.Range("I122").Value = 1787.5834
$ws.Range("K122").Value = 5362.7502
$ws.Range("M122").Value = -2912.7502
$ws.Range("H126").Value = 1956.2
$ws.Range("I126").Value = 1509.174
$ws.Range("K126").Value = 4527.522
$ws.Range("M126").Value = -2057.522
$ws.Range("H136").Value = 4010
$ws.Range("I136").Value = 814.625
$ws.Range("K136").Value = 2443.875
$ws.Range("M136").Value = 106.125
